$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width (column B) ---
$ws.Columns.Item(2).ColumnWidth = 12.43

# --- Re-style rows 2 & 3 (title/subtitle) and B6:C6 label to border+center ---
$ws.Range("B2:D3").HorizontalAlignment = -4108
$ws.Range("B6:C6").HorizontalAlignment = -4108

# --- Reorder header row 4: Recuperados, Fallecidos, Enfermos ---
$ws.Range("B4").Value = "Recuperados"
$ws.Range("C4").Value = "Fallecidos"
$ws.Range("D4").Value = "Enfermos"

# --- Selection ---
$ws.Range("D5").Select()
